$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.5224089492919877
$ws.Range("D3").Value = 0.05346793379164971
$ws.Range("D4").Value = 0.1928148357833853
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0.1622197797262909
$ws.Range("D6").Value = 0.1164530646275216
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 0.7189786359051346
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0.3405075779509109
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0.4035442946966274
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 0.5875736605387843
$ws.Range("D11").Value = 0.7539191170245518
$ws.Range("D12").Value = 0.6724653949994213
$ws.Range("D13").Value = 0.4130369763631875
$ws.Range("D14").Value = 0.8566938781689506
$ws.Range("D15").Value = 0.5350776907514587
$ws.Range("D16").Value = 0.2194665035703925
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 0.2839788186421354
$ws.Range("D18").Value = 0.4214600720225786
$ws.Range("D19").Value = 0.4726430697686028
$ws.Range("D20").Value = 0.1753210991973831
$ws.Range("D21").Value = 0.3430851417457589
$ws.Range("D22").Value = 0.6814860075287087
$ws.Range("D23").Value = 0.1529615361683906
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 0.2796500698685587
$ws.Range("D25").Value = 0.1751927649632402
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 0.285062567305792
$ws.Range("D27").Value = 0.4769332991020206
$ws.Range("D28").Value = 0.1248631480793064
$ws.Range("D29").Value = 0.5191615702080846
$ws.Range("D30").Value = 0.344626024746535
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = 0.4484590674112255
$ws.Range("D32").Value = 0.742536086640561
$ws.Range("D33").Value = 0.1404297445956508
$ws.Range("D34").Value = 0.1517662380094652
$ws.Range("D35").Value = 0.8679707707472764
$ws.Range("D36").Value = 0.2981665168296947
$ws.Range("D37").Value = 0.7272128687317856
$ws.Range("D38").Value = 0.7452467764189785
$ws.Range("D39").Value = 0.1789677327648747
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 0.4005172546844303
$ws.Range("C41").Value = 1
$ws.Range("D41").Value = 0.6731047751413147
$ws.Range("D42").Value = 0.5228099320666804
$ws.Range("D43").Value = 0.8717630636676866
$ws.Range("D44").Value = 0.7894012744015866
$ws.Range("C45").Value = 1
$ws.Range("D45").Value = 0.5296866409084169
$ws.Range("D46").Value = 0.8024705896301483
$ws.Range("D47").Value = 0.2868161318068595
$ws.Range("C48").Value = 1
$ws.Range("D48").Value = 0.5863925254753707
$ws.Range("D49").Value = 0.9062319982473641
